$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking values (e.g. "156.25")
# are not auto-converted to numbers by Excel type inference, matching the
# original inline/shared string cell type. Style is reset to Normal afterwards
# so no stray numFmt/style index lingers on these cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.494.33'
$ws.Range("E2").Value = '  +1.66%  '

$ws.Range("D3").Value = '2.288.25'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '156.25'
$ws.Range("E5").Value = '  +15,507.75%  '

$ws.Range("D6").Value = '307.35'
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").Value = '95.92'
$ws.Range("E7").Value = '  +4.45%  '

$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '0.496'
$ws.Range("E10").Value = '  +2.87%  '

$ws.Range("D11").Value = '35.75'
$ws.Range("E11").Value = '  +10.11%  '

$ws.Range("D12").Value = '0.0806'
$ws.Range("E12").Value = '  +1.31%  '

$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("E14").Value = '  +2.18%  '

$ws.Range("D15").Value = '2.643.78'
$ws.Range("E15").Value = '  +1.10%  '

$ws.Range("D16").Value = '14.53'
$ws.Range("E16").Value = '  +2.38%  '

$ws.Range("D17").Value = '2.291.52'
$ws.Range("E17").Value = '  +0.80%  '

$ws.Range("D18").Value = '0.803'
$ws.Range("E18").Value = '  +5.43%  '

$ws.Range("D19").Value = '42.362.37'
$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("E20").Value = '  +4.08%  '

$ws.Range("E21").Value = '  +2.01%  '

$ws.Range("E22").Value = '  +2.14%  '

$ws.Range("D23").Value = '68.17'
$ws.Range("E23").Value = '  +2.07%  '

$ws.Range("D24").Value = '243.20'
$ws.Range("E24").Value = '  +1.03%  '

$ws.Range("D25").Value = '2.62'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("E26").Value = '  +2.20%  '

$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").Value = '24.15'
$ws.Range("E28").Value = '  +0.05%  '

$ws.Range("D29").Value = '36.46'
$ws.Range("E29").Value = '  +7.81%  '

$ws.Range("D30").Value = '9.65'
$ws.Range("E30").Value = '  +0.84%  '

$ws.Range("E31").Value = '  -8.46%  '

$ws.Range("D32").Value = '161.37'
$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("E33").Value = '  +4.20%  '

$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = '0.0755'
$ws.Range("E35").Value = '  +1.38%  '

$ws.Range("D36").Value = '3.09'
$ws.Range("E36").Value = '  +2.84%  '

$ws.Range("E37").Value = '  +5.27%  '

$ws.Range("D38").Value = '17.20'
$ws.Range("E38").Value = '  +1.80%  '

$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  +1.73%  '

$ws.Range("D42").Value = '4.20'
$ws.Range("E42").Value = '  +7.59%  '

$ws.Range("D43").Value = '2.020.25'
$ws.Range("E43").Value = '  -1.97%  '

$ws.Range("D44").Value = '19.75'
$ws.Range("E44").Value = '  +1.70%  '

$ws.Range("D46").Value = '0.0284'
$ws.Range("E46").Value = '  +2.11%  '

$ws.Range("D47").Value = '10.23'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("E48").Value = '  +3.81%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '1.55'
$ws.Range("E49").Value = '  +1.55%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '53.44'
$ws.Range("E50").Value = '  +3.62%  '

$ws.Range("D51").Value = '73.43'
$ws.Range("E51").Value = '  +1.23%  '

# Reset number format back to the default (General) style so the written
# cells keep no explicit style index, same as the rest of the sheet.
$ws.Range("D2:E51").Style = "Normal"
